$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q6)
$ws.Range("B7").Value = 0.1643430888929492
$ws.Range("C7").Value = 2.713245774742223
$ws.Range("D7").Value = 10.87156620932192
$ws.Range("E7").Value = 3.297205818465374
$ws.Range("F7").Value = 3.337312348085441
$ws.Range("G7").Value = 38

# Row 8 (Q7)
$ws.Range("B8").Value = 0.3953194985244485
$ws.Range("C8").Value = 2.739053368442637
$ws.Range("D8").Value = 12.07707068967368
$ws.Range("E8").Value = 3.475208006677252
$ws.Range("F8").Value = 3.500275178861427
$ws.Range("G8").Value = 37

# Row 9 (Q8)
$ws.Range("B9").Value = 0.2568635966504935
$ws.Range("C9").Value = 1.810696804862654
$ws.Range("D9").Value = 5.367596856534075
$ws.Range("E9").Value = 2.316807470752388
$ws.Range("F9").Value = 2.362340041758294
$ws.Range("G9").Value = 20

# Row 10 (Q9)
$ws.Range("B10").Value = -0.4680654340455334
$ws.Range("C10").Value = 1.381782322298316
$ws.Range("D10").Value = 2.829787416625887
$ws.Range("E10").Value = 1.682197199089895
$ws.Range("F10").Value = 1.681743345435759
$ws.Range("G10").Value = 13

# Row 11 (Q9 duplicate label)
$ws.Range("B11").Value = -1.213277339139011
$ws.Range("C11").Value = 2.013103686783858
$ws.Range("D11").Value = 5.583994774231593
$ws.Range("E11").Value = 2.363047772312611
$ws.Range("F11").Value = 2.267143817825457
$ws.Range("G11").Value = 5
